# Update data collection model:
# - rename prefixed header columns to their short form on the
#   BiomedicalConcept and DataElementConcept sheets
# - drop the duplicated data-validation rules that were left over
#   on those same sheets
# - remove the accidental duplicate sheets (BiomedicalConcept1, Coding1,
#   DataElementConcept1)
# - make sure the first sheet is the active one again

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- BiomedicalConcept sheet -------------------------------------------------
$ws1 = $wb.Worksheets.Item("BiomedicalConcept")
$ws1.Range("C1").Value = "conceptId"
$ws1.Range("D1").Value = "ncitCode"
$ws1.Range("E1").Value = "href"

# The sheet had each validation rule duplicated; remove them and re-add a
# single copy of each so the sheet ends up with exactly one "B" rule and
# one "J" rule.
$ws1.Range("B2:B1048576").Validation.Delete()
$ws1.Range("J2:J1048576").Validation.Delete()

$bcList = $ws1.Range("B2:B1048576").Validation
$bcList.Add(3, 1, 1, '"bc"')
$bcList.ShowInput = $false
$bcList.ShowError = $false

$resultScales = $ws1.Range("J2:J1048576").Validation
$resultScales.Add(3, 1, 1, '"Ordinal,Narrative,Nominal,Quantitative,Temporal"')
$resultScales.ShowInput = $false
$resultScales.ShowError = $false

# --- DataElementConcept sheet ------------------------------------------------
$ws3 = $wb.Worksheets.Item("DataElementConcept")
$ws3.Range("A1").Value = "conceptId"
$ws3.Range("B1").Value = "ncitCode"
$ws3.Range("C1").Value = "href"

# Same deduplication for the single validation rule on this sheet.
$ws3.Range("E2:E1048576").Validation.Delete()
$dataType = $ws3.Range("E2:E1048576").Validation
$dataType.Add(3, 1, 1, '"boolean,date,datetime,decimal,duration,integer,string,uri"')
$dataType.ShowInput = $false
$dataType.ShowError = $false

# --- Remove the duplicate "1" sheets -----------------------------------------
$wb.Worksheets.Item("BiomedicalConcept1").Delete()
$wb.Worksheets.Item("Coding1").Delete()
$wb.Worksheets.Item("DataElementConcept1").Delete()

# --- Restore the active tab to the first sheet -------------------------------
$wb.Worksheets.Item("BiomedicalConcept").Activate()
